$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getRelationById")

# Insert a new column C ("label" / relation type filter) shifting the
# existing relationId/rspCode/rspMessage/rspStatus columns one to the right.
$ws.Columns("C").Insert()

# Header + sample data for the newly inserted "label" column.
$ws.Range("C1").Value = "label"
$ws.Range("C2").Value = "invert"
$ws.Range("C3").Value = "Has_Device"

# Replace the old threaded comment (originally anchored on C1, now logically
# on D1 after the column insert) with updated wording, and add a new
# threaded comment explaining the new "label" column on C1.
$ws.Range("C1").Comment.Delete()

$ws.Range("C1").AddCommentThreaded('If this field is set, test case will use its value to get a list of relations, then pick up the id of the 1st relation as the relationId value to be test and the input parameter "relationId" will be ignored.') | Out-Null

$ws.Range("D1").AddCommentThreaded("Note: relation ids are different in dev/test/prod environment, so use a fixed id here may lead to test fail in a another environment!") | Out-Null

# Move the active sheet / selection to getRelationById (was updateEntity).
$ws.Activate() | Out-Null
$ws.Range("D14").Select() | Out-Null
